$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that must stay text (not auto-converted to numbers),
# so force text format on those cells before assigning.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.959.96"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.382.86"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.93"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.98"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.53"
$ws.Range("E10").Value = "  -6.27%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.15"
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.745.82"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.376.64"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.818"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.69"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.877.94"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0956"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.19"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.35"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  -4.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.51"
$ws.Range("E27").Value = "  -9.69%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.80"
$ws.Range("E30").Value = "  +20.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.08"
$ws.Range("E31").Value = "  +4.48%  "
$ws.Range("E32").Value = "  +6.88%  "
$ws.Range("E33").Value = "  -4.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.13"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0773"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  +6.79%  "
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("E39").Value = "  -5.02%  "
$ws.Range("E40").Value = "  -5.11%  "
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("E42").Value = "  -7.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.939.23"
$ws.Range("E43").Value = "  +4.77%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.48"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("E46").Value = "  -10.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.47"
$ws.Range("E47").Value = "  +5.51%  "
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.51"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.616.73"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "68.50"
$ws.Range("E51").Value = "  -8.34%  "
